# Change representation of missing values in Survey_Data
# Empty cells that represent "missing" values get the literal text "-"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey_Data")

# List of cells that were empty and now hold the string "-"
$cellsToDash = @(
    "E2", "F2",
    "E4", "F4",
    "E5", "F5",
    "E28", "F28",
    "E32", "F32",
    "E34", "F34",
    "E39", "F39",
    "A40",
    "E40", "F40",
    "E46", "F46"
)

foreach ($addr in $cellsToDash) {
    $ws.Range($addr).Value = "-"
}

# F46 also picks up the same (unaccented) font style used by its row neighbors
# (matching the style shift seen from s="1" to s="2" in the saved XML)
$ws.Range("E46").Copy()
$ws.Range("F46").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F46").Value = "-"
$excel.CutCopyMode = $false

# Update the view so the selection matches the final file
# (the workbook also scrolls so row 38 is the top-left visible row;
# the COM surface here doesn't expose a scroll-position setter, so only
# the selection - which is the part reflected in the cell-level state - is set)
$ws.Activate()
$ws.Range("D41").Select()
